{"js": "// Office.js (Word JavaScript API) implementation of the edit described\n// by the diff:\n//   1. Ronish Barahi's roll number changes from \"000\" to \"067\".\n//   2. Nirajan Sha's roll number changes from \"000\" to \"054\".\n//   3. A \"_GoBack\" bookmark is added to the (empty) paragraph that comes\n//      right after the \"DATE: 2079-04-24\" paragraph.\n//\n// (The remaining hunks in the diff are either re-save artifacts of the\n// authoring tool -- e.g. pruned/added namespace declarations on the root\n// elements of document.xml/styles.xml/numbering.xml/endnotes.xml/\n// footnotes.xml, dropped latentStyle \"Mention/Smart Hyperlink/...\"\n// exceptions, a dropped w16cid:durableId / w15:restartNumberingAfterBreak\n// attribute -- or a run split in the title paragraph that keeps the\n// visible text 100% identical. None of those are observable/achievable\n// through the Word object model, so they are intentionally not\n// reproduced here.)\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Ronish Barahi: (078 BCT 000) -> (078 BCT 067)\n// ---------------------------------------------------------------------\nconst ronishContext = body.search(\"Ronish Barahi (078 BCT 000)\", { matchCase: true });\nronishContext.load(\"items\");\nawait context.sync();\n\nif (ronishContext.items.length === 0) {\n  throw new Error(\"Could not find 'Ronish Barahi (078 BCT 000)' in the document.\");\n}\n\nconst ronishNumber = ronishContext.items[0].search(\"000\", { matchCase: true });\nronishNumber.load(\"items\");\nawait context.sync();\nronishNumber.items[0].insertText(\"067\", Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Nirajan Sha: (078 BCT 000) -> (078 BCT 054)\n// ---------------------------------------------------------------------\nconst nirajanContext = body.search(\"Nirajan Sha (078 BCT 000)\", { matchCase: true });\nnirajanContext.load(\"items\");\nawait context.sync();\n\nif (nirajanContext.items.length === 0) {\n  throw new Error(\"Could not find 'Nirajan Sha (078 BCT 000)' in the document.\");\n}\n\nconst nirajanNumber = nirajanContext.items[0].search(\"000\", { matchCase: true });\nnirajanNumber.load(\"items\");\nawait context.sync();\nnirajanNumber.items[0].insertText(\"054\", Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) Insert the \"_GoBack\" bookmark into the empty paragraph that follows\n//    the \"DATE: 2079-04-24\" paragraph.\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet dateParagraphIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"DATE:\") !== -1) {\n    dateParagraphIndex = i;\n    break;\n  }\n}\n\nif (dateParagraphIndex === -1 || dateParagraphIndex + 1 >= paragraphs.items.length) {\n  throw new Error(\"Could not locate the paragraph following the DATE line.\");\n}\n\nconst bookmarkParagraph = paragraphs.items[dateParagraphIndex + 1];\nconst bookmarkRange = bookmarkParagraph.getRange();\nbookmarkRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) implementation of the edit described by the\n# diff:\n#   1. Ronish Barahi's roll number changes from \"000\" to \"067\".\n#   2. Nirajan Sha's roll number changes from \"000\" to \"054\".\n#   3. A \"_GoBack\" bookmark is added to the (empty) paragraph that comes\n#      right after the \"DATE: 2079-04-24\" paragraph.\n#\n# (The remaining hunks in the diff are either re-save artifacts of the\n# authoring tool -- e.g. pruned/added namespace declarations on the root\n# elements of document.xml/styles.xml/numbering.xml/endnotes.xml/\n# footnotes.xml, dropped latentStyle \"Mention/Smart Hyperlink/...\"\n# exceptions, a dropped w16cid:durableId / w15:restartNumberingAfterBreak\n# attribute -- or a run split in the title paragraph that keeps the\n# visible text 100% identical. None of those are observable/achievable\n# through the Word object model, so they are intentionally not\n# reproduced here.)\n\n$d = $word.ActiveDocument\n$paragraphs = @($d.Paragraphs)\n\n# ---------------------------------------------------------------------\n# 1) Ronish Barahi: (078 BCT 000) -> (078 BCT 067)\n# ---------------------------------------------------------------------\n$ronishParagraph = $null\nforeach ($p in $paragraphs) {\n    if ($p.Range.Text.Contains(\"Ronish Barahi\") -and $p.Range.Text.Contains(\"078 BCT 000\")) {\n        $ronishParagraph = $p\n        break\n    }\n}\nif ($null -eq $ronishParagraph) {\n    throw \"Could not find the 'Ronish Barahi (078 BCT 000)' paragraph.\"\n}\n$ronishFind = $ronishParagraph.Range.Find\n$ronishFind.ClearFormatting()\n$ronishFind.Text = \"000\"\n$ronishFind.Replacement.ClearFormatting()\n$ronishFind.Replacement.Text = \"067\"\n$ronishFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 2) Nirajan Sha: (078 BCT 000) -> (078 BCT 054)\n# ---------------------------------------------------------------------\n$nirajanParagraph = $null\nforeach ($p in $paragraphs) {\n    if ($p.Range.Text.Contains(\"Nirajan Sha\") -and $p.Range.Text.Contains(\"078 BCT 000\")) {\n        $nirajanParagraph = $p\n        break\n    }\n}\nif ($null -eq $nirajanParagraph) {\n    throw \"Could not find the 'Nirajan Sha (078 BCT 000)' paragraph.\"\n}\n$nirajanFind = $nirajanParagraph.Range.Find\n$nirajanFind.ClearFormatting()\n$nirajanFind.Text = \"000\"\n$nirajanFind.Replacement.ClearFormatting()\n$nirajanFind.Replacement.Text = \"054\"\n$nirajanFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3) Insert the \"_GoBack\" bookmark into the empty paragraph that follows\n#    the \"DATE: 2079-04-24\" paragraph.\n# ---------------------------------------------------------------------\n$dateParagraphIndex = -1\nfor ($i = 0; $i -lt $paragraphs.Count; $i++) {\n    if ($paragraphs[$i].Range.Text.Contains(\"DATE:\")) {\n        $dateParagraphIndex = $i\n        break\n    }\n}\nif ($dateParagraphIndex -eq -1 -or ($dateParagraphIndex + 1) -ge $paragraphs.Count) {\n    throw \"Could not locate the paragraph following the DATE line.\"\n}\n$bookmarkParagraph = $paragraphs[$dateParagraphIndex + 1]\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkParagraph.Range) | Out-Null\n"}
